$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.5797639999999999
$ws.Range("M2").Value = 71.67520133333333
$ws.Range("N2").Value = 215.025604
$ws.Range("O2").Value = 0.7843164235321155
$ws.Range("P2").Value = 0.7843164235321155
$ws.Range("Q2").Value = 13.85156714193955
$ws.Range("R2").Value = 124.664104277456
$ws.Range("S2").Value = 0.7843164235321155
$ws.Range("T2").Value = 0.7843164235321155

# Row 3
$ws.Range("H3").Value = 0.5797639999999999
$ws.Range("O3").Value = 0.1685558319854606
$ws.Range("P3").Value = 0.1685558319854606
$ws.Range("S3").Value = 0.1685558319854606
$ws.Range("T3").Value = 0.1685558319854606

# Row 4
$ws.Range("H4").Value = 0.5797639999999999
$ws.Range("M4").Value = 1.213835666666667
$ws.Range("N4").Value = 3.641507
$ws.Range("O4").Value = 0.0132825751602454
$ws.Range("P4").Value = 0.0132825751602454
$ws.Range("Q4").Value = 0.2345794071497777
$ws.Range("R4").Value = 2.111214664348
$ws.Range("S4").Value = 0.0132825751602454
$ws.Range("T4").Value = 0.0132825751602454

# Row 5
$ws.Range("H5").Value = 0.5797639999999999
$ws.Range("M5").Value = 1.940298333333333
$ws.Range("N5").Value = 5.820895
$ws.Range("O5").Value = 0.02123199964668382
$ws.Range("P5").Value = 0.02123199964668382
$ws.Range("Q5").Value = 0.3749717076422222
$ws.Range("R5").Value = 3.37474536878
$ws.Range("S5").Value = 0.02123199964668382
$ws.Range("T5").Value = 0.02123199964668382

# Row 6
$ws.Range("H6").Value = 0.5797639999999999
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.081443333333333
$ws.Range("N6").Value = 3.24433
$ws.Range("O6").Value = 0.01183385259719093
$ws.Range("P6").Value = 0.01183385259719093
$ws.Range("Q6").Value = 0.2089939709022222
$ws.Range("R6").Value = 1.88094573812
$ws.Range("S6").Value = 0.01183385259719093
$ws.Range("T6").Value = 0.01183385259719093

# Row 7
$ws.Range("H7").Value = 0.5797639999999999
$ws.Range("M7").Value = 0.07121833333333334
$ws.Range("N7").Value = 0.213655
$ws.Range("O7").Value = 0.0007793170783036338
$ws.Range("P7").Value = 0.0007793170783036338
$ws.Range("Q7").Value = 0.01376327526888889
$ws.Range("R7").Value = 0.12386947742
$ws.Range("S7").Value = 0.0007793170783036338
$ws.Range("T7").Value = 0.0007793170783036338
